$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

$ws.Range("A18").Value = "22/02/2018"
$ws.Range("B18").Value = "3070"
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = "anvil test2"
$ws.Range("E18").Value = "1"
$ws.Range("F18").Value = "80000571"
$ws.Range("G18").Value = "7"

$ws.Range("A19").Value = "22/02/2018"
$ws.Range("B19").Value = "3070"
$ws.Range("C19").Value = 99
$ws.Range("D19").Value = "anvil test2"
$ws.Range("E19").Value = "1"
$ws.Range("F19").Value = "80000571"
$ws.Range("G19").Value = "8"

$ws.Range("A20").Value = "22/02/2018"
$ws.Range("B20").Value = "3070"
$ws.Range("C20").Value = 100
$ws.Range("D20").Value = "anvil test2"
$ws.Range("E20").Value = "1"
$ws.Range("F20").Value = "80000571"
$ws.Range("G20").Value = "9"

$ws.Range("A21").Value = "22/02/2018"
$ws.Range("B21").Value = "3067"
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = "anvil test"
$ws.Range("E21").Value = "1"
$ws.Range("F21").Value = "80000571"
$ws.Range("G21").Value = "16"

$ws.Range("A22").Value = "22/02/2018"
$ws.Range("B22").Value = "3067"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = "anvil test2"
$ws.Range("E22").Value = "1"
$ws.Range("F22").Value = "80000571"
$ws.Range("G22").Value = "16"

$ws.Range("A23").Value = "22/02/2018"
$ws.Range("B23").Value = "3067"
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = "anvil test"
$ws.Range("E23").Value = "1"
$ws.Range("F23").Value = "80000571"
$ws.Range("G23").Value = "17"

$ws.Range("A24").Value = "22/02/2018"
$ws.Range("B24").Value = "3067"
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = "anvil test2"
$ws.Range("E24").Value = "1"
$ws.Range("F24").Value = "80000571"
$ws.Range("G24").Value = "17"

$ws.Range("A25").Value = "22/02/2018"
$ws.Range("B25").Value = "3067"
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = "anvil test"
$ws.Range("E25").Value = "1"
$ws.Range("F25").Value = "80000571"
$ws.Range("G25").Value = "18"

$ws.Range("A26").Value = "**"
